$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 75: new data row appended after running the profit script on 2026-02-07
# Column A holds a date-like string; force it to stay plain text (not get
# auto-converted to a date serial number) by temporarily using a text number
# format, then reset the style back to Normal so no style index lingers on
# the cell (matching the rest of the date column).
$cellA = $ws.Cells.Item(75, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "02/07/2026"
$cellA.Style = "Normal"

$ws.Cells.Item(75, 2).Value = 9279.24
$ws.Cells.Item(75, 3).Value = 0.2383866365551371
$ws.Cells.Item(75, 4).Value = 0.7616133634448629
$ws.Cells.Item(75, 5).Value = -319.06
$ws.Cells.Item(75, 6).Value = -38.21
$ws.Cells.Item(75, 7).Value = -23824.56
$ws.Cells.Item(75, 8).Value = -77.12
$ws.Cells.Item(75, 9).Value = -1090.8
$ws.Cells.Item(75, 10).Value = -33.03
$ws.Cells.Item(75, 11).Value = -24915.36
$ws.Cells.Item(75, 12).Value = -72.86
